$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new bullet paragraph ("CI/CD Tool: Jenkins") right before the
#    "Repository management tool" bullet.
# ---------------------------------------------------------------------------
$target = $d.Content
$found = $target.Find.Execute("Repository management tool", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Repository management tool' paragraph"
}

$repoPara = $target.Paragraphs(1)

# Creates a new empty paragraph right before $repoPara, duplicating its pPr
# (ListParagraph style, numbering, tabs, spacing).
$target.InsertParagraphBefore()

# Re-locate "Repository management tool" (its index shifted by one) and grab
# the newly-created empty paragraph immediately preceding it.
$target2 = $d.Content
$found2 = $target2.Find.Execute("Repository management tool", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$repoIdx = $target2.Paragraphs(1).Index
$newPara = $d.Paragraphs($repoIdx - 1)
$newRange = $newPara.Range

$ciCdXml = @"
<?xml version='1.0'?>
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val='ListParagraph'/>
<w:numPr><w:ilvl w:val='0'/><w:numId w:val='6'/></w:numPr>
<w:tabs><w:tab w:val='left' w:pos='583'/></w:tabs>
<w:spacing w:before='17'/>
<w:rPr><w:sz w:val='18'/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:color w:val='5883AF'/><w:sz w:val='18'/></w:rPr><w:t>CI/CD Tool</w:t></w:r>
<w:r><w:rPr><w:color w:val='5883AF'/><w:sz w:val='18'/></w:rPr><w:t>:</w:t></w:r>
<w:r><w:rPr><w:sz w:val='18'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>
<w:r><w:rPr><w:color w:val='5883AF'/><w:sz w:val='18'/></w:rPr><w:t>Jenkins</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$newRange.InsertXML($ciCdXml)

# ---------------------------------------------------------------------------
# 2) Fix "BitBucket" -> "Bitbucket" in the "Repository management tool"
#    bullet, and drop the now-stale spell-check <w:proofErr/> markers around
#    it.
# ---------------------------------------------------------------------------
$target3 = $d.Content
$found3 = $target3.Find.Execute("Repository management tool", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not re-find 'Repository management tool' paragraph"
}

$repoPara2 = $target3.Paragraphs(1)
$repoRange = $repoPara2.Range

$repoXml = @"
<?xml version='1.0'?>
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'>
<w:body>
<w:p w14:paraId='78118427' w14:textId='6E8BA6F8' w:rsidR='00090B5D' w:rsidRPr='00090B5D' w:rsidRDefault='00090B5D'>
<w:pPr>
<w:pStyle w:val='ListParagraph'/>
<w:numPr><w:ilvl w:val='0'/><w:numId w:val='6'/></w:numPr>
<w:tabs><w:tab w:val='left' w:pos='583'/></w:tabs>
<w:spacing w:before='17'/>
<w:rPr><w:sz w:val='18'/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:b/><w:color w:val='5883AF'/><w:sz w:val='18'/></w:rPr><w:t>Repository management tool</w:t></w:r>
<w:r w:rsidRPr='00090B5D'><w:rPr><w:sz w:val='18'/></w:rPr><w:t>:</w:t></w:r>
<w:r><w:rPr><w:sz w:val='18'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>
<w:r w:rsidRPr='00090B5D'><w:rPr><w:color w:val='5883AF'/><w:sz w:val='18'/></w:rPr><w:t>Bitbucket</w:t></w:r>
<w:r w:rsidR='00AC45C1'><w:rPr><w:color w:val='5883AF'/><w:sz w:val='18'/></w:rPr><w:t>, GIT</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$repoRange.InsertXML($repoXml)
